$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (Price / Volume(1h) columns) per upstream refresh.
# Column D (Price) holds plain-text numeric-looking strings in the source file
# (e.g. thousand-dot-grouped "71.809.30", which is not a valid number anyway and
# stays text automatically -- but plain decimals like "596.13" would otherwise be
# auto-coerced to a real number by the COM value setter). For those we briefly
# force the cell to Text format, assign the value, then restore the default style
# so the saved cell keeps its original (unstyled) text representation, matching
# the source.

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $ws.Range($cell.Address()).Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = '71.809.30'
$ws.Cells.Item(2, 5).Value = '  +0.94%  '
$ws.Cells.Item(3, 4).Value = '2.685.58'
$ws.Cells.Item(3, 5).Value = '  +2.56%  '
Set-TextValue 5 4 '596.13'
$ws.Cells.Item(5, 5).Value = '  -1.47%  '
Set-TextValue 6 4 '175.04'
$ws.Cells.Item(6, 5).Value = '  -2.24%  '
$ws.Cells.Item(7, 5).Value = '  +0.02%  '
$ws.Cells.Item(8, 5).Value = '  -0.60%  '
$ws.Cells.Item(9, 4).Value = '2.685.17'
$ws.Cells.Item(9, 5).Value = '  +2.62%  '
$ws.Cells.Item(10, 5).Value = '  +2.03%  '
$ws.Cells.Item(11, 5).Value = '  +2.39%  '
Set-TextValue 12 4 '0.352'
$ws.Cells.Item(12, 5).Value = '  +1.44%  '
$ws.Cells.Item(13, 5).Value = '  -1.16%  '
$ws.Cells.Item(14, 4).Value = '3.157.24'
$ws.Cells.Item(14, 5).Value = '  +2.10%  '
$ws.Cells.Item(15, 5).Value = '  -0.95%  '
$ws.Cells.Item(16, 4).Value = '71.797.33'
$ws.Cells.Item(16, 5).Value = '  +0.97%  '
Set-TextValue 17 4 '26.18'
$ws.Cells.Item(17, 5).Value = '  -1.76%  '
$ws.Cells.Item(18, 4).Value = '2.639.40'
$ws.Cells.Item(18, 5).Value = '  +0.36%  '
Set-TextValue 19 4 '12.08'
$ws.Cells.Item(19, 5).Value = '  +5.22%  '
Set-TextValue 20 4 '8.02'
$ws.Cells.Item(20, 5).Value = '  +2.05%  '
Set-TextValue 21 4 '370.67'
$ws.Cells.Item(21, 5).Value = '  -2.73%  '
$ws.Cells.Item(22, 5).Value = '  -0.36%  '
Set-TextValue 23 4 '2.02'
$ws.Cells.Item(23, 5).Value = '  +1.27%  '
Set-TextValue 24 4 '71.73'
$ws.Cells.Item(24, 5).Value = '  -0.79%  '
$ws.Cells.Item(25, 5).Value = '  +0.01%  '
Set-TextValue 26 4 '4.30'
$ws.Cells.Item(26, 5).Value = '  -3.37%  '
Set-TextValue 27 4 '9.82'
$ws.Cells.Item(27, 5).Value = '  -1.56%  '
$ws.Cells.Item(28, 4).Value = '2.721.75'
$ws.Cells.Item(28, 5).Value = '  -1.16%  '
Set-TextValue 29 4 '0.999'
$ws.Cells.Item(29, 5).Value = '  -0.17%  '
$ws.Cells.Item(30, 4).Value = '0.0₃0942'
$ws.Cells.Item(30, 5).Value = '  -1.20%  '
Set-TextValue 31 4 '8.01'
$ws.Cells.Item(31, 5).Value = '  -0.31%  '
Set-TextValue 32 4 '504.04'
$ws.Cells.Item(32, 5).Value = '  -7.90%  '
$ws.Cells.Item(33, 5).Value = '  -3.61%  '
Set-TextValue 34 4 '1.81'
$ws.Cells.Item(34, 5).Value = '  -0.97%  '
$ws.Cells.Item(35, 5).Value = '  -0.03%  '
Set-TextValue 36 4 '163.71'
$ws.Cells.Item(36, 5).Value = '  -0.96%  '
Set-TextValue 37 4 '19.46'
$ws.Cells.Item(37, 5).Value = '  +1.38%  '
$ws.Cells.Item(38, 5).Value = '  -0.41%  '
Set-TextValue 39 4 '1.37'
$ws.Cells.Item(39, 5).Value = '  -2.55%  '
Set-TextValue 40 4 '0.107'
$ws.Cells.Item(40, 5).Value = '  -6.47%  '
Set-TextValue 41 4 '1.79'
$ws.Cells.Item(41, 5).Value = '  -4.24%  '
$ws.Cells.Item(42, 5).Value = '  -0.03%  '
Set-TextValue 43 4 '4.98'
$ws.Cells.Item(43, 5).Value = '  -0.60%  '
Set-TextValue 44 4 '2.55'
$ws.Cells.Item(44, 5).Value = '  -2.45%  '
Set-TextValue 45 4 '0.331'
$ws.Cells.Item(45, 5).Value = '  -0.27%  '
Set-TextValue 46 4 '39.16'
$ws.Cells.Item(46, 5).Value = '  -1.99%  '
Set-TextValue 47 4 '154.31'
$ws.Cells.Item(47, 5).Value = '  +1.01%  '
Set-TextValue 48 4 '3.71'
$ws.Cells.Item(48, 5).Value = '  +1.96%  '
Set-TextValue 49 4 '0.547'
$ws.Cells.Item(49, 5).Value = '  +2.57%  '
Set-TextValue 50 4 '1.73'
$ws.Cells.Item(50, 5).Value = '  +3.13%  '
Set-TextValue 51 4 '0.0764'
$ws.Cells.Item(51, 5).Value = '  +1.20%  '
